# Update header/label row (row 2) to the new field names, clearing the
# trailing two columns that no longer carry a header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "person name"
$ws.Range("B2").Value = "phone number"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""

# Row 4's "man" flag was being written as an empty string instead of a real
# boolean FALSE - fix it to set the proper zero/false value.
$ws.Range("D4").Value = $false
